$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1973.1428
$ws.Range("I98").Value = 1973.1428
$ws.Range("K98").Value = 1973.1428
$ws.Range("M98").Value = -475.1428000000001
$ws.Range("H103").Value = 395.2857
$ws.Range("I103").Value = 316.75
$ws.Range("K103").Value = 950.25
$ws.Range("M103").Value = -364.25
$ws.Range("H112").Value = 3712.5925
$ws.Range("J112").Value = 3717.6
$ws.Range("L112").Value = 11152.8
$ws.Range("N112").Value = -13368.8
$ws.Range("H122").Value = 1973.1428
$ws.Range("I122").Value = 1973.1428
$ws.Range("K122").Value = 5919.428400000001
$ws.Range("M122").Value = -3469.428400000001
$ws.Range("H132").Value = 8961.5
$ws.Range("I132").Value = 11147.625
$ws.Range("J132").Value = 2403.125
$ws.Range("K132").Value = 33442.875
$ws.Range("L132").Value = 7209.375
$ws.Range("M132").Value = -30912.875
$ws.Range("N132").Value = -12269.375
$ws.Range("H141").Value = 3448.5
$ws.Range("I141").Value = 2252.6365
$ws.Range("J141").Value = 7833.3335
$ws.Range("K141").Value = 6757.9095
$ws.Range("L141").Value = 23500.0005
$ws.Range("M141").Value = -1577.9095
$ws.Range("N141").Value = -33860.00049999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 10092.857
$ws.Range("I31").Value = 7608.5
$ws.Range("K31").Value = 7608.5
$ws.Range("M31").Value = -7314.5
$ws.Range("H61").Value = 4768813.5
$ws.Range("I61").Value = 9123.091
$ws.Range("J61").Value = 10004473
$ws.Range("K61").Value = 9123.091
$ws.Range("L61").Value = 10004473
$ws.Range("M61").Value = -8911.091
$ws.Range("N61").Value = -10004897
$ws.Range("H132").Value = 387037.5
$ws.Range("I132").Value = 448786.47
$ws.Range("K132").Value = 1346359.41
$ws.Range("M132").Value = -1343829.41
$ws.Range("H136").Value = 4768813.5
$ws.Range("I136").Value = 9123.091
$ws.Range("J136").Value = 10004473
$ws.Range("K136").Value = 27369.273
$ws.Range("L136").Value = 30013419
$ws.Range("M136").Value = -24819.273
$ws.Range("N136").Value = -30018519

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3151163
$ws.Range("I134").Value = 4721.4443
$ws.Range("K134").Value = 14164.3329
$ws.Range("M134").Value = -11629.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5565870
$ws.Range("I58").Value = 4208.5386
$ws.Range("J58").Value = 41716668
$ws.Range("K58").Value = 4208.5386
$ws.Range("L58").Value = 41716668
$ws.Range("M58").Value = -4005.5386
$ws.Range("N58").Value = -41717074
$ws.Range("H94").Value = 8154.5
$ws.Range("I94").Value = 34766.332
$ws.Range("J94").Value = 2013.3077
$ws.Range("K94").Value = 34766.332
$ws.Range("L94").Value = 2013.3077
$ws.Range("M94").Value = -34315.332
$ws.Range("N94").Value = -2915.3077
$ws.Range("H99").Value = 50327.6
$ws.Range("J99").Value = 40888.5
$ws.Range("L99").Value = 40888.5
$ws.Range("N99").Value = -43884.5
$ws.Range("H126").Value = 50327.6
$ws.Range("J126").Value = 40888.5
$ws.Range("L126").Value = 122665.5
$ws.Range("N126").Value = -127605.5
$ws.Range("H132").Value = 2641.9312
$ws.Range("I132").Value = 2672.238
$ws.Range("K132").Value = 8016.714
$ws.Range("M132").Value = -5486.714
$ws.Range("H136").Value = 5565870
$ws.Range("I136").Value = 4208.5386
$ws.Range("J136").Value = 41716668
$ws.Range("K136").Value = 12625.6158
$ws.Range("L136").Value = 125150004
$ws.Range("M136").Value = -10075.6158
$ws.Range("N136").Value = -125155104

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2337.5
$ws.Range("I39").Value = 3
$ws.Range("J39").Value = 2517.077
$ws.Range("K39").Value = 9
$ws.Range("L39").Value = 7551.231000000001
$ws.Range("M39").Value = 285
$ws.Range("N39").Value = -8139.231000000001
$ws.Range("H137").Value = 5186.8213
$ws.Range("I137").Value = 3247.6667
$ws.Range("J137").Value = 7424.3076
$ws.Range("K137").Value = 9743.000100000001
$ws.Range("L137").Value = 22272.9228
$ws.Range("M137").Value = -4643.000100000001
$ws.Range("N137").Value = -32472.9228
$ws.Range("H139").Value = 4910
$ws.Range("I139").Value = 3551.9092
$ws.Range("K139").Value = 10655.7276
$ws.Range("M139").Value = -5515.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1581.1177
$ws.Range("I102").Value = 1575.4814
$ws.Range("K102").Value = 1575.4814
$ws.Range("M102").Value = 46.51860000000011
$ws.Range("H122").Value = 4326.852
$ws.Range("I122").Value = 4712.353
$ws.Range("J122").Value = 3671.5
$ws.Range("K122").Value = 14137.059
$ws.Range("L122").Value = 11014.5
$ws.Range("M122").Value = -11687.059
$ws.Range("N122").Value = -15914.5
$ws.Range("H132").Value = 10779.596
$ws.Range("I132").Value = 9069.108
$ws.Range("K132").Value = 27207.324
$ws.Range("M132").Value = -24677.324
$ws.Range("H136").Value = 81173.086
$ws.Range("J136").Value = 81173.086
$ws.Range("L136").Value = 243519.258
$ws.Range("N136").Value = -248619.258

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7023.9165
$ws.Range("I7").Value = 10300
$ws.Range("K7").Value = 10300
$ws.Range("M7").Value = -10188
$ws.Range("H93").Value = 2178.3333
$ws.Range("I93").Value = 1780.8667
$ws.Range("K93").Value = 1780.8667
$ws.Range("M93").Value = -532.8667
$ws.Range("H122").Value = 2872
$ws.Range("I122").Value = 2595.318
$ws.Range("K122").Value = 7785.954000000001
$ws.Range("M122").Value = -5335.954000000001
$ws.Range("H126").Value = 7023.9165
$ws.Range("I126").Value = 10300
$ws.Range("K126").Value = 30900
$ws.Range("M126").Value = -28430
$ws.Range("H132").Value = 4871382.5
$ws.Range("I132").Value = 6493799
$ws.Range("K132").Value = 19481397
$ws.Range("M132").Value = -19478867
$ws.Range("H134").Value = 60136
$ws.Range("J134").Value = 60136
$ws.Range("L134").Value = 60136
$ws.Range("N134").Value = -70276
$ws.Range("H136").Value = 5687731
$ws.Range("I136").Value = 4633077
$ws.Range("K136").Value = 13899231
$ws.Range("M136").Value = -13896681

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 4277.4
$ws.Range("I107").Value = 5260.25
$ws.Range("J107").Value = 346
$ws.Range("K107").Value = 15780.75
$ws.Range("L107").Value = 1038
$ws.Range("M107").Value = -13860.75
$ws.Range("N107").Value = -4878
$ws.Range("H132").Value = 4506536
$ws.Range("I132").Value = 6411964.5
$ws.Range("J132").Value = 2795.182
$ws.Range("K132").Value = 19235893.5
$ws.Range("L132").Value = 8385.545999999998
$ws.Range("M132").Value = -19233363.5
$ws.Range("N132").Value = -13445.546

Write-Output "Applied all profit updates"